# NIT-9009066343.xlsx — "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The underlying edit corrects the "Periodo Mora" (arrears period) / "Valor Mora"
# (arrears amount) pairs recorded for two workers, and nudges the company logo
# slightly to the left.
#
# Worker GILBERTO ANDRES ZABALETA VILORIA (rows 16-18) had periods 1707/1706/1705;
# the corrected data swaps the 1707 and 1705 rows (with their matching arrears
# amounts) so the periods read 1705/1706/1707 in order, while row 17 (1706) is
# unaffected.
#
# Worker CARLOS ARMANDO BUESACO DIAZ GRANADOS (rows 20-21) had periods 2010/2009;
# the corrected data swaps them to read 2009/2010 in order (the arrears amount is
# identical for both periods, so only the period text changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- GILBERTO ANDRES ZABALETA VILORIA: swap row 16 <-> row 18 (Periodo Mora / Valor Mora)
$ws.Range("E16").Value = "1705"
$ws.Range("F16").Value = 98933
$ws.Range("E18").Value = "1707"
$ws.Range("F18").Value = 106000

# --- CARLOS ARMANDO BUESACO DIAZ GRANADOS: swap row 20 <-> row 21 (Periodo Mora)
$ws.Range("E20").Value = "2009"
$ws.Range("E21").Value = "2010"

# --- Nudge the company logo 13.5pt to the left (same size, same vertical position)
$logo = $ws.Shapes.Item(1)
$logo.Left = $logo.Left - 13.5
